$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data window: 30.07.2024 09:00 through 06.08.2024 09:00 (169 hourly rows)
$colA = @(45503,45503,45503,45503,45503,45503,45503,45503,45503,45503,45503,45503,45503,45503,45503,45504,45504,45504,45504,45504,45504,45504,45504,45504,45504,45504,45504,45504,45504,45504,45504,45504,45504,45504,45504,45504,45504,45504,45504,45505,45505,45505,45505,45505,45505,45505,45505,45505,45505,45505,45505,45505,45505,45505,45505,45505,45505,45505,45505,45505,45505,45505,45505,45506,45506,45506,45506,45506,45506,45506,45506,45506,45506,45506,45506,45506,45506,45506,45506,45506,45506,45506,45506,45506,45506,45506,45506,45507,45507,45507,45507,45507,45507,45507,45507,45507,45507,45507,45507,45507,45507,45507,45507,45507,45507,45507,45507,45507,45507,45507,45507,45508,45508,45508,45508,45508,45508,45508,45508,45508,45508,45508,45508,45508,45508,45508,45508,45508,45508,45508,45508,45508,45508,45508,45508,45509,45509,45509,45509,45509,45509,45509,45509,45509,45509,45509,45509,45509,45509,45509,45509,45509,45509,45509,45509,45509,45509,45509,45509,45510,45510,45510,45510,45510,45510,45510,45510,45510,45510)
$colB = @(9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,0,1,2,3,4,5,6,7,8,9)
$colC = @(0,0,2.798,3.088,3.267,3.161,2.94,2.444,1.868,1.203,0.466,0.065,0,0,0,0,0,0,0,0,0,0.066,0.368,1.017,1.756,2.404,2.747,2.936,3.139,3.035,2.893,2.345,1.805,1.166,0.452,0.063,0,0,0,0,0,0,0,0,0,0.046,0.368,0.955,1.59,2.335,2.437,2.915,3.084,3.085,2.911,2.33,1.723,1.02,0.425,0.045,0,0,0,0,0,0,0,0,0,0.04,0.247,0.53,0.978,1.48,1.753,2.061,2.176,2.126,2.027,1.836,1.422,0.823,0.292,0.036,0,0,0,0,0,0,0,0,0,0.046,0.364,0.895,1.486,1.625,1.884,1.878,1.824,1.767,1.637,1.52,1.126,0.789,0.292,0.038,0,0,0,0,0,0,0,0,0,0.039,0.319,0.656,1.525,2.129,2.61,2.788,2.933,2.84,2.715,2.139,1.602,1.004,0.367,0.041,0,0,0,0,0,0,0,0,0,0.039,0.306,0.788,1.532,2.208,2.594,2.873,3.086,2.841,2.579,2.057,1.533,0.896,0.364,0.041,0,0,0,0,0,0,0,0,0,0.042,0.333,0.957,1.628)
$colD = @("30.07.20249","30.07.202410","30.07.202411","30.07.202412","30.07.202413","30.07.202414","30.07.202415","30.07.202416","30.07.202417","30.07.202418","30.07.202419","30.07.202420","30.07.202421","30.07.202422","30.07.202423","31.07.20240","31.07.20241","31.07.20242","31.07.20243","31.07.20244","31.07.20245","31.07.20246","31.07.20247","31.07.20248","31.07.20249","31.07.202410","31.07.202411","31.07.202412","31.07.202413","31.07.202414","31.07.202415","31.07.202416","31.07.202417","31.07.202418","31.07.202419","31.07.202420","31.07.202421","31.07.202422","31.07.202423","01.08.20240","01.08.20241","01.08.20242","01.08.20243","01.08.20244","01.08.20245","01.08.20246","01.08.20247","01.08.20248","01.08.20249","01.08.202410","01.08.202411","01.08.202412","01.08.202413","01.08.202414","01.08.202415","01.08.202416","01.08.202417","01.08.202418","01.08.202419","01.08.202420","01.08.202421","01.08.202422","01.08.202423","02.08.20240","02.08.20241","02.08.20242","02.08.20243","02.08.20244","02.08.20245","02.08.20246","02.08.20247","02.08.20248","02.08.20249","02.08.202410","02.08.202411","02.08.202412","02.08.202413","02.08.202414","02.08.202415","02.08.202416","02.08.202417","02.08.202418","02.08.202419","02.08.202420","02.08.202421","02.08.202422","02.08.202423","03.08.20240","03.08.20241","03.08.20242","03.08.20243","03.08.20244","03.08.20245","03.08.20246","03.08.20247","03.08.20248","03.08.20249","03.08.202410","03.08.202411","03.08.202412","03.08.202413","03.08.202414","03.08.202415","03.08.202416","03.08.202417","03.08.202418","03.08.202419","03.08.202420","03.08.202421","03.08.202422","03.08.202423","04.08.20240","04.08.20241","04.08.20242","04.08.20243","04.08.20244","04.08.20245","04.08.20246","04.08.20247","04.08.20248","04.08.20249","04.08.202410","04.08.202411","04.08.202412","04.08.202413","04.08.202414","04.08.202415","04.08.202416","04.08.202417","04.08.202418","04.08.202419","04.08.202420","04.08.202421","04.08.202422","04.08.202423","05.08.20240","05.08.20241","05.08.20242","05.08.20243","05.08.20244","05.08.20245","05.08.20246","05.08.20247","05.08.20248","05.08.20249","05.08.202410","05.08.202411","05.08.202412","05.08.202413","05.08.202414","05.08.202415","05.08.202416","05.08.202417","05.08.202418","05.08.202419","05.08.202420","05.08.202421","05.08.202422","05.08.202423","06.08.20240","06.08.20241","06.08.20242","06.08.20243","06.08.20244","06.08.20245","06.08.20246","06.08.20247","06.08.20248","06.08.20249")

$startRow = 2
for ($i = 0; $i -lt $colA.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $colA[$i]
    $ws.Cells.Item($r, 2).Value = $colB[$i]
    $ws.Cells.Item($r, 3).Value = $colC[$i]
    $ws.Cells.Item($r, 4).Value = $colD[$i]
}
